# Replace "OBSERVER" terminology with "SUBSCRIBER" terminology in the
# responselist worksheet (commit: "changed observer message to subscriber message")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C15").Value = "SUBSCRIBER REGISTERED"
$ws.Range("C16").Value = "SUBSCRIBER UNREGISTERED"
$ws.Range("C28").Value = "SUBSCRIBER ALREADY REGISTERED"
$ws.Range("C29").Value = "SUBSCRIBER NOT REGISTERED"
$ws.Range("C31").Value = "SUBSCRIBER MISSING REGESTRATION STATE"
$ws.Range("C32").Value = "SUBSCRIBER CANT REGISTER AS DEVICE"
$ws.Range("C33").Value = "DEVICE CANT REGISTER AS SUBSCRIBER"
